$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from H1 (bold, bordered, centered) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data rows 2-22 for columns I (I0) and J (IF)
$iValues = @{
    2  = 2
    3  = 1
    4  = 4
    5  = 7
    6  = 4
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
}

$jValues = @{
    2  = 4
    3  = 4
    4  = 6
    5  = 8
    6  = 6
    7  = 6
    8  = 3
    9  = 5
    10 = 4
    11 = 4
    12 = 6
    13 = 6
    14 = 5
    15 = 6
    16 = 5
    17 = 4
    18 = 6
    19 = 7
    20 = 7
    21 = 3
    22 = 2
}

foreach ($r in 2..22) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}
